$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# "REMISIONES   ENERO  2022  " (January 2022) is worksheet #4
$wsEnero = $sheets.Item(4)
# "REMISIONES FEBRERO   2022  " (February 2022) is worksheet #5 (the originally active sheet)
$wsFebrero = $sheets.Item(5)

# ---------------------------------------------------------------------------
# Sheet "REMISIONES   ENERO  2022  " (sheet4): row 45 gets marked as paid
# ---------------------------------------------------------------------------
$wsEnero.Range("F45").Value = 44608
$wsEnero.Range("G45").Value = 18072

# ---------------------------------------------------------------------------
# Sheet "REMISIONES FEBRERO   2022  " (sheet5): mark rows 21 and 24 as paid,
# and fill in the previously-blank rows 33-49 with the February closing data.
# ---------------------------------------------------------------------------
$wsFebrero.Range("F21").Value = 44611
$wsFebrero.Range("G21").Value = 48706

$wsFebrero.Range("F24").Value = 44608
$wsFebrero.Range("G24").Value = 2655

$wsFebrero.Range("A33").Value = 44608
$wsFebrero.Range("D33").Value = "PROSUBCA"
$wsFebrero.Range("E33").Value = 1520
$wsFebrero.Range("F33").Value = 44615
$wsFebrero.Range("G33").Value = 1520

$wsFebrero.Range("A34").Value = 44608
$wsFebrero.Range("D34").Value = "OBRADOR"
$wsFebrero.Range("E34").Value = 7496

$wsFebrero.Range("A35").Value = 44608
$wsFebrero.Range("D35").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E35").Value = 1060

$wsFebrero.Range("A36").Value = 44609
$wsFebrero.Range("D36").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E36").Value = 6970

$wsFebrero.Range("A37").Value = 44610
$wsFebrero.Range("D37").Value = "OBRADOR"
$wsFebrero.Range("E37").Value = 16986
$wsFebrero.Range("F37").Value = 44610
$wsFebrero.Range("G37").Value = 16986

$wsFebrero.Range("A38").Value = 44610
$wsFebrero.Range("D38").Value = "ISRAEL LEDO"
$wsFebrero.Range("E38").Value = 40856

$wsFebrero.Range("A39").Value = 44610
$wsFebrero.Range("D39").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E39").Value = 7105

$wsFebrero.Range("A40").Value = 44611
$wsFebrero.Range("D40").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E40").Value = 7690

$wsFebrero.Range("A41").Value = 44611
$wsFebrero.Range("D41").Value = "OBRADOR"
$wsFebrero.Range("E41").Value = 340

$wsFebrero.Range("A42").Value = 44613
$wsFebrero.Range("D42").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E42").Value = 8456

$wsFebrero.Range("A43").Value = 44613
$wsFebrero.Range("D43").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E43").Value = 21362

$wsFebrero.Range("A44").Value = 44614
$wsFebrero.Range("D44").Value = "OBRADOR"
$wsFebrero.Range("E44").Value = 11560

$wsFebrero.Range("A45").Value = 44615
$wsFebrero.Range("D45").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E45").Value = 189945

$wsFebrero.Range("A46").Value = 44615
$wsFebrero.Range("D46").Value = "PROSUBCA"
$wsFebrero.Range("E46").Value = 1313

$wsFebrero.Range("A47").Value = 44615
$wsFebrero.Range("D47").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E47").Value = 91036

$wsFebrero.Range("A48").Value = 44615
$wsFebrero.Range("D48").Value = "ABASTOS  HERRADURA "
$wsFebrero.Range("E48").Value = 170

$wsFebrero.Range("A49").Value = 44615
$wsFebrero.Range("D49").Value = "COMERCIO CENTRAL "
$wsFebrero.Range("E49").Value = 400

# Recalculate all dependent formulas (SUM totals, H column shared formulas, etc.)
$excel.Calculate()

# ---------------------------------------------------------------------------
# Restore view/selection state to match the closing session: the user ended
# up reviewing cell G46 on the January sheet before returning to the
# February sheet (which stays the active tab) with D50 selected.
# ---------------------------------------------------------------------------
$wsEnero.Activate()
$wsEnero.Range("G46").Select()

$wsFebrero.Activate()
$wsFebrero.Range("D50").Select()
